$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 953.55554
$ws.Range("I98").Value = 947.875
$ws.Range("K98").Value = 947.875
$ws.Range("M98").Value = 550.125
$ws.Range("H122").Value = 953.55554
$ws.Range("I122").Value = 947.875
$ws.Range("K122").Value = 2843.625
$ws.Range("M122").Value = -393.625
$ws.Range("H132").Value = 1672.1428
$ws.Range("I132").Value = 1705.9
$ws.Range("K132").Value = 5117.700000000001
$ws.Range("M132").Value = -2587.700000000001
$ws.Range("H137").Value = 998.2
$ws.Range("I137").Value = 998.25
$ws.Range("K137").Value = 2994.75
$ws.Range("M137").Value = -444.75
$ws.Range("H141").Value = 2895.818
$ws.Range("I141").Value = 1910
$ws.Range("J141").Value = 4621
$ws.Range("K141").Value = 5730
$ws.Range("L141").Value = 13863
$ws.Range("M141").Value = -550
$ws.Range("N141").Value = -24223

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1683
$ws.Range("I5").Value = 1979.6
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 1979.6
$ws.Range("L5").Value = 200
$ws.Range("M5").Value = -1867.6
$ws.Range("N5").Value = -424
$ws.Range("H32").Value = 3350988.8
$ws.Range("I32").Value = 3336749.5
$ws.Range("K32").Value = 3336749.5
$ws.Range("M32").Value = -3336462.5
$ws.Range("H50").Value = 5754.8335
$ws.Range("I50").Value = 8361.5
$ws.Range("J50").Value = 541.5
$ws.Range("K50").Value = 8361.5
$ws.Range("L50").Value = 541.5
$ws.Range("M50").Value = -7647.5
$ws.Range("N50").Value = -1969.5
$ws.Range("H97").Value = 845.6667
$ws.Range("I97").Value = 882.6429000000001
$ws.Range("J97").Value = 716.25
$ws.Range("K97").Value = 882.6429000000001
$ws.Range("L97").Value = 716.25
$ws.Range("M97").Value = -386.6429000000001
$ws.Range("N97").Value = -1708.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1683
$ws.Range("I4").Value = 1979.6
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 1979.6
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = -1864.6
$ws.Range("N4").Value = -430
$ws.Range("H107").Value = 2707.318
$ws.Range("I107").Value = 2517.2354
$ws.Range("J107").Value = 3353.6
$ws.Range("K107").Value = 2517.2354
$ws.Range("L107").Value = 3353.6
$ws.Range("M107").Value = -597.2354
$ws.Range("N107").Value = -7193.6
$ws.Range("H134").Value = 5948.5713
$ws.Range("I134").Value = 5768
$ws.Range("J134").Value = 6400
$ws.Range("K134").Value = 17304
$ws.Range("L134").Value = 19200
$ws.Range("M134").Value = -14769
$ws.Range("N134").Value = -24270

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 1181.5
$ws.Range("J11").Value = 1472.25
$ws.Range("L11").Value = 1472.25
$ws.Range("N11").Value = -1752.25
$ws.Range("H16").Value = 5503.3335
$ws.Range("I16").Value = 4508
$ws.Range("K16").Value = 4508
$ws.Range("M16").Value = -4221
$ws.Range("H31").Value = 1344.875
$ws.Range("I31").Value = 1189.75
$ws.Range("J31").Value = 1500
$ws.Range("K31").Value = 1189.75
$ws.Range("L31").Value = 1500
$ws.Range("M31").Value = -894.75
$ws.Range("N31").Value = -2090
$ws.Range("H34").Value = 1344.875
$ws.Range("I34").Value = 1189.75
$ws.Range("J34").Value = 1500
$ws.Range("K34").Value = 1189.75
$ws.Range("L34").Value = 1500
$ws.Range("M34").Value = -987.75
$ws.Range("N34").Value = -1904
$ws.Range("H86").Value = 5981.778
$ws.Range("I86").Value = 6179.5
$ws.Range("K86").Value = 6179.5
$ws.Range("M86").Value = -5056.5
$ws.Range("H89").Value = 5981.778
$ws.Range("I89").Value = 6179.5
$ws.Range("K89").Value = 30897.5
$ws.Range("M89").Value = -25281.5
$ws.Range("H113").Value = 5503.3335
$ws.Range("I113").Value = 4508
$ws.Range("K113").Value = 4508
$ws.Range("M113").Value = -2338

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8636315
$ws.Range("I4").Value = 10001618
$ws.Range("J4").Value = 3175103.2
$ws.Range("K4").Value = 30004854
$ws.Range("L4").Value = 9525309.600000001
$ws.Range("M4").Value = -30004742
$ws.Range("N4").Value = -9525533.600000001
$ws.Range("H23").Value = 1003.3
$ws.Range("I23").Value = 776.6667
$ws.Range("J23").Value = 1100.4286
$ws.Range("K23").Value = 2330.0001
$ws.Range("L23").Value = 3301.2858
$ws.Range("M23").Value = -2095.0001
$ws.Range("N23").Value = -3771.2858
$ws.Range("H56").Value = 11232.8125
$ws.Range("I56").Value = 11232.8125
$ws.Range("K56").Value = 11232.8125
$ws.Range("M56").Value = -10702.8125
$ws.Range("H75").Value = 2299.3333
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("H78").Value = 2299.3333
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("N78").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8001.3335
$ws.Range("I70").Value = 6999.5
$ws.Range("K70").Value = 6999.5
$ws.Range("M70").Value = -6729.5
$ws.Range("H73").Value = 8001.3335
$ws.Range("I73").Value = 6999.5
$ws.Range("K73").Value = 6999.5
$ws.Range("M73").Value = -6063.5
$ws.Range("H93").Value = 20251
$ws.Range("J93").Value = 20251
$ws.Range("L93").Value = 20251
$ws.Range("N93").Value = -23995
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("H132").Value = 3334.875
$ws.Range("J132").Value = 3371.5
$ws.Range("L132").Value = 10114.5
$ws.Range("N132").Value = -15174.5
$ws.Range("N100").ClearContents()
$ws.Range("N106").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 873.875
$ws.Range("J22").Value = 973.5
$ws.Range("L22").Value = 973.5
$ws.Range("N22").Value = -1563.5
$ws.Range("H27").Value = 873.875
$ws.Range("J27").Value = 973.5
$ws.Range("L27").Value = 973.5
$ws.Range("N27").Value = -1187.5
$ws.Range("H46").Value = 1838.7
$ws.Range("I46").Value = 1398.4375
$ws.Range("J46").Value = 3599.75
$ws.Range("K46").Value = 1398.4375
$ws.Range("L46").Value = 3599.75
$ws.Range("M46").Value = -1210.4375
$ws.Range("N46").Value = -3975.75
$ws.Range("H93").Value = 1014.65216
$ws.Range("I93").Value = 1182.4667
$ws.Range("K93").Value = 1182.4667
$ws.Range("M93").Value = 65.53330000000005
$ws.Range("H100").Value = 4538
$ws.Range("I100").Value = 4860.4
$ws.Range("J100").Value = 4000.6667
$ws.Range("K100").Value = 4860.4
$ws.Range("L100").Value = 4000.6667
$ws.Range("M100").Value = -4319.4
$ws.Range("N100").Value = -5082.6667
$ws.Range("H122").Value = 6518.5586
$ws.Range("I122").Value = 5123.2
$ws.Range("J122").Value = 7620.1577
$ws.Range("K122").Value = 15369.6
$ws.Range("L122").Value = 22860.4731
$ws.Range("M122").Value = -12919.6
$ws.Range("N122").Value = -27760.4731
$ws.Range("H136").Value = 1937.6923
$ws.Range("I136").Value = 1838
$ws.Range("K136").Value = 5514
$ws.Range("M136").Value = -2964

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 20194.5
$ws.Range("J104").Value = 20194.5
$ws.Range("L104").Value = 20194.5
$ws.Range("N104").Value = -27182.5
$ws.Range("H139").Value = 73216.336
$ws.Range("I139").Value = 79650
$ws.Range("J139").Value = 69999.5
$ws.Range("K139").Value = 79650
$ws.Range("L139").Value = 69999.5
$ws.Range("M139").Value = -74510
$ws.Range("N139").Value = -80279.5
